# Edit script: add "2022-Q3" sheet data (feat: add 2022-Q3 data)
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Update the "总计" (totals) summary sheet: insert the 2022-Q3 row at the
#    top of the data (row 2) and shift the existing quarters down.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")

# Row 9 is brand new (the sheet previously only had rows 1-8); give its A
# cell the same "index column" formatting used by A2:A8 before filling it in.
$totalWs.Cells.Item(8, 1).Copy()
$totalWs.Cells.Item(9, 1).PasteSpecial(-4122)

$totalsData = @(
    @(0, '2022-Q3', 22, 4.86),
    @(1, '2022-Q2', 28, 9.640000000000001),
    @(2, '2022-Q1', 47, 9.18),
    @(3, '2021-Q4', 43, 6.83),
    @(4, '2021-Q3', 55, 14.14),
    @(5, '2021-Q2', 157, 26.77),
    @(6, '2021-Q1', 50, 11.19),
    @(7, '2020-Q4', 35, 10.32)
)

foreach ($row in $totalsData) {
    $r = [int]$row[0] + 2
    $totalWs.Cells.Item($r, 1).Value = [int]$row[0]
    $totalWs.Cells.Item($r, 2).Value = $row[1]
    $totalWs.Cells.Item($r, 3).Value = $row[2]
    $totalWs.Cells.Item($r, 4).Value = $row[3]
}

# ---------------------------------------------------------------------------
# 2. Insert the new "2022-Q3" worksheet right after "总计" (i.e. before the
#    existing "2022-Q2" sheet). Copying the "2022-Q2" sheet gives us the
#    exact same layout/styles (header row + index column formatting) as the
#    other quarterly sheets, then we overwrite its contents with 2022-Q3 data.
# ---------------------------------------------------------------------------
$sourceWs = $wb.Worksheets.Item("2022-Q2")
$afterWs = $wb.Worksheets.Item("总计")
$sourceWs.Copy($null, $afterWs)

$newWs = $wb.Worksheets.Item(2)
$newWs.Name = "2022-Q3"

$fundData = @(
    @(0, '000628', '大成高新技术产业股票A', '43.11', '75.37', '3.24', '1.3968', 10),
    @(1, '008269', '大成睿享混合A', '19.80', '66.91', '4.34', '0.8593', 4),
    @(2, '010892', '中银证券精选行业股票A', '10.84', '93.07', '5.83', '0.6320', 7),
    @(3, '090007', '大成策略回报混合', '9.90', '61.86', '6.32', '0.6257', 1),
    @(4, '090013', '大成竞争优势混合', '6.88', '61.00', '4.07', '0.2800', 4),
    @(5, '011066', '大成高新技术产业股票C', '6.69', '75.37', '3.24', '0.2168', 10),
    @(6, '011834', '大成投资严选六月持有混合A', '3.10', '66.75', '6.68', '0.2071', 2),
    @(7, '013463', '大成致远优势一年持有期混合A', '3.65', '60.88', '5.48', '0.2000', 3),
    @(8, '008270', '大成睿享混合C', '4.02', '66.91', '4.34', '0.1745', 4),
    @(9, '009640', '中银证券优选行业龙头混合A', '1.01', '93.07', '8.23', '0.0831', 5),
    @(10, '004495', '博时量化平衡混合', '4.26', '39.09', '0.74', '0.0315', 10),
    @(11, '009246', '大摩ESG量化混合', '2.69', '84.29', '1.00', '0.0269', 7),
    @(12, '009641', '中银证券优选行业龙头混合C', '0.32', '93.07', '8.23', '0.0263', 5),
    @(13, '005055', '华泰柏瑞量化阿尔法灵活配置混合A', '1.98', '92.07', '1.20', '0.0238', 8),
    @(14, '010893', '中银证券精选行业股票C', '0.40', '93.07', '5.83', '0.0233', 7),
    @(15, '011835', '大成投资严选六月持有混合C', '0.22', '66.75', '6.68', '0.0147', 2),
    @(16, '005444', '光大保德信多策略精选18个月定期开放灵活配置混合', '0.84', '29.28', '1.63', '0.0137', 5),
    @(17, '013464', '大成致远优势一年持有期混合C', '0.17', '60.88', '5.48', '0.0093', 3),
    @(18, '001397', '建信精工制造指数增强', '0.46', '89.20', '2.03', '0.0093', 6),
    @(19, '501069', '华宝标普中国Ａ股质量价值指数（LOF）', '0.14', '93.78', '3.67', '0.0051', 4),
    @(20, '001900', '诺安精选价值混合', '0.12', '85.83', '2.66', '0.0032', 8),
    @(21, '006532', '华泰柏瑞量化阿尔法灵活配置混合C', '0.01', '92.07', '1.20', '0.0001', 8)
)

# Force columns B (fund code, has leading zeros) and D:G (figures stored as
# text in the source data) to stay text instead of auto-converting to numbers.
$newWs.Range("B2:B23").NumberFormat = "@"
$newWs.Range("D2:G23").NumberFormat = "@"

foreach ($row in $fundData) {
    $r = [int]$row[0] + 2
    $newWs.Cells.Item($r, 1).Value = [int]$row[0]
    $newWs.Cells.Item($r, 2).Value = $row[1]
    $newWs.Cells.Item($r, 3).Value = $row[2]
    $newWs.Cells.Item($r, 4).Value = $row[3]
    $newWs.Cells.Item($r, 5).Value = $row[4]
    $newWs.Cells.Item($r, 6).Value = $row[5]
    $newWs.Cells.Item($r, 7).Value = $row[6]
    $newWs.Cells.Item($r, 8).Value = [int]$row[7]
}

# The copied sheet had 28 data rows (2022-Q2); 2022-Q3 only has 22, so clear
# the leftover rows (24-29) to shrink the used range back down to A1:H23.
$newWs.Range("A24:H29").Clear()

# Copying a sheet makes the new copy the active tab; restore the workbook's
# original active-tab position (3rd sheet, 1-based) now that "2022-Q2" sits
# there instead of "2022-Q1".
$wb.Worksheets.Item(3).Activate()

Write-Output "done"
